$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 43; this shifts the existing rows 43-49 down to 44-50,
# carrying all their existing values/formatting with them.
$ws.Rows("43:43").Insert()

# --- Fill in the brand-new row 43 (weekly update entry) ---
$ws.Range("A43").Value = 3
$ws.Range("B43").Value = "Femacal de La Calera"
$ws.Range("C43").Value = "Coquimbo"
$ws.Range("D43").Value = 44889
$ws.Range("E43").Value = 5
$ws.Range("F43").Value = 300000000
$ws.Range("G43").Value = "Espárragos"
$ws.Range("H43").Value = "Verde"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 1330
$ws.Range("K43").Value = 1300
$ws.Range("L43").Value = 1400
$ws.Range("M43").Value = 1349
$ws.Range("N43").Value = "$/kilo"
$ws.Range("O43").Value = "Provincia de Quillota"
$ws.Range("P43").Value = 1349
$ws.Range("Q43").Value = 1
$ws.Range("R43").Value = "Hortaliza"

# --- Update the dates that shifted on the rows pushed down by the insert ---
$ws.Range("D44").Value = 44889
$ws.Range("D46").Value = 44161
$ws.Range("D48").Value = 44172
